# Europe_3_Factors.xlsx -- "Final Update: project is accomplished"
#
# Column A currently holds numeric YYYYMM period codes (e.g. 201201).
# Convert each one into a text label formatted as "YYYY-MM" (e.g. "2012-01"),
# derived straight from the existing stored value so every row (2-93) is
# handled without hard-coding the date sequence.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 93; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    $raw = [string]$cell.Value2
    $yyyy = $raw.Substring(0, 4)
    $mm = $raw.Substring(4, 2)
    $cell.Value = "$yyyy-$mm"
}

# Move the active selection to I41 (matches the saved view state in the diff).
$ws.Range("I41").Select()
